# Changes for Auth filter
# - Reorders / renames existing sheets
# - Inserts a new "login" sheet (register -> login -> getUsers -> getInvitations -> sendInvitations)
# - Populates the login sheet with request/response sample data

$wb = $excel.ActiveWorkbook

# --- Rename + reorder existing sheets -------------------------------------
# Original order: GetAllAppUsers, SaveOrUpdateUSer, GetAllInvitation, SendInvitation
# Target order:   register, login, getUsers, getInvitations, sendInvitations

# Move "SaveOrUpdateUSer" to the front of the tab strip, then rename it.
$wb.Worksheets.Item("SaveOrUpdateUSer").Move($wb.Worksheets.Item("GetAllAppUsers")) | Out-Null
$wb.Worksheets.Item("SaveOrUpdateUSer").Name = "register"

$wb.Worksheets.Item("GetAllAppUsers").Name = "getUsers"
$wb.Worksheets.Item("GetAllInvitation").Name = "getInvitations"
$wb.Worksheets.Item("SendInvitation").Name = "sendInvitations"

# --- Insert the new "login" sheet right after "register" ------------------
$loginSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("register"))
$loginSheet.Name = "login"

# request / response sample payloads (kept CRLF like the rest of the workbook)
$loginRequest  = "{`r`n    ""username"" : ""gsk@admin.com"",`r`n    ""password"" : ""assword""`r`n}"
$loginResponse = "{`r`n    ""token"": ""eyJhbGciOiJIUzI1NiJ9.eyJzdWIiOiJnc2tAYWRtaW4uY29tIiwiZXhwIjoxNjAzOTA1NDQwLCJpYXQiOjE2MDM4Njk0NDB9.TeC02BXSXPfAY45sj4p9xzW_GpwXLXWg5q0EynCoVNs"",`r`n    ""message"": ""Login success"",`r`n    ""results"": null`r`n}"

$loginSheet.Range("A1").Value = "login"
$loginSheet.Range("B1").Value = "POST"
# Set D1 (response) before C1 (request) so the shared-string table is built
# in the same order as the target workbook (response before request).
$loginSheet.Range("D1").Value = $loginResponse
$loginSheet.Range("C1").Value = $loginRequest
$loginSheet.Range("C1:D1").WrapText = $true
$loginSheet.Rows.Item(1).RowHeight = 144

# --- Cosmetic touch-ups (selection / active tab) ---------------------------
$ws = $wb.Worksheets.Item("register")
$ws.Activate() | Out-Null
$ws.Range("C8").Select() | Out-Null

$ws = $loginSheet
$ws.Activate() | Out-Null
$ws.Range("F7").Select() | Out-Null

$ws = $wb.Worksheets.Item("getUsers")
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null

$ws = $wb.Worksheets.Item("sendInvitations")
$ws.Activate() | Out-Null
$ws.Range("G17").Select() | Out-Null

Write-Host "Auth filter changes applied"
